# DTSCCI-248 use external short name rather than court venue name
#
# The merge-field placeholder:
#   <<caseManagementLocation.venue_name>><<else>> Online Civil Claims<<es_>>
# becomes:
#   <<caseManagementLocation.external_short_name>><<else>>Online Civil Claims<<es_>>
#
# i.e. "venue_name" -> "external_short_name", and the stray space that used
# to separate "<<else>>" from "Online Civil Claims" is removed.

$d = $word.ActiveDocument

$old = ".venue_name>><<else>> Online Civil Claims<<es_>>"
$new = ".external_short_name>><<else>>Online Civil Claims<<es_>>"

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find the venue_name placeholder text to replace"
}
